$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'AA221452PDFN"
$ws.Range("E2").Value = "TT23200FW369 12:0"
$ws.Range("D2").Value = "PASSED"
$ws.Range("F2").Value = "19 jul. 2023, 12:08:30"

$ws.Columns.Item(1).ColumnWidth = 7.1796875
$ws.Columns.Item(3).ColumnWidth = 14.90625

$ws.Range("F10").Select()
